# Update cryptocurrency price/volume data in the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D='51.681.07'; E='  +4.66%  ' }
    3 = @{ D='2.761.52'; E='  +4.64%  ' }
    4 = @{ E='  -0.14%  ' }
    5 = @{ D='117.76'; E='  +5.77%  ' }
    6 = @{ D='333.05'; E='  +3.06%  ' }
    7 = @{ D='0.535'; E='  +2.63%  ' }
    8 = @{ E='  -0.06%  ' }
    9 = @{ D='0.576'; E='  +6.14%  ' }
    10 = @{ D='41.50'; E='  +4.14%  ' }
    11 = @{ E='  +1.75%  ' }
    12 = @{ E='  +2.50%  ' }
    13 = @{ E='  +3.13%  ' }
    14 = @{ E='  +5.44%  ' }
    15 = @{ D='3.192.33'; E='  +4.26%  ' }
    16 = @{ D='2.767.00'; E='  +3.94%  ' }
    17 = @{ D='0.887'; E='  +2.95%  ' }
    18 = @{ D='51.623.30'; E='  +4.47%  ' }
    19 = @{ D='13.71'; E='  +6.63%  ' }
    20 = @{ D='2.99'; E='  +3.14%  ' }
    21 = @{ D='6.87'; E='  +3.04%  ' }
    22 = @{ D='0.0₃0965'; E='  +2.36%  ' }
    23 = @{ D='277.35'; E='  +1.83%  ' }
    24 = @{ D='70.21'; E='  +0.01%  ' }
    25 = @{ D='2.67'; E='  +5.32%  ' }
    26 = @{ D='26.92'; E='  +2.20%  ' }
    27 = @{ D='4.15'; E='  +0.42%  ' }
    28 = @{ D='0.999'; E='  +0.08%  ' }
    29 = @{ D='10.31'; E='  +2.20%  ' }
    30 = @{ E='  -0.47%  ' }
    31 = @{ E='  +2.40%  ' }
    32 = @{ D='35.64'; E='  +0.84%  ' }
    33 = @{ D='50.50'; E='  +1.96%  ' }
    34 = @{ D='5.61'; E='  +2.91%  ' }
    35 = @{ D='0.0837'; E='  +5.45%  ' }
    36 = @{ D='19.43'; E='  +1.40%  ' }
    37 = @{ E='  +4.28%  ' }
    38 = @{ E='  -0.29%  ' }
    39 = @{ D='5.00'; E='  +0.32%  ' }
    40 = @{ D='3.27'; E='  +5.53%  ' }
    41 = @{ D='130.64'; E='  +4.47%  ' }
    42 = @{ D='23.54'; E='  +7.17%  ' }
    43 = @{ D='0.0345'; E='  +10.16%  ' }
    44 = @{ E='  +2.94%  ' }
    45 = @{ D='2.28'; E='  +3.37%  ' }
    46 = @{ D='2.38'; E='  +11.57%  ' }
    47 = @{ D='2.116.56'; E='  +0.92%  ' }
    48 = @{ D='3.36'; E='  +3.78%  ' }
    49 = @{ D='2.28'; E='  +3.35%  ' }
    50 = @{ D='5.63'; E='  +7.70%  ' }
    51 = @{ D='9.01'; E='  +1.16%  ' }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    if ($vals.ContainsKey("D")) {
        # Force text format first so numeric-looking strings (e.g. "2.99", "0.999")
        # are not auto-converted to numbers by Excel when assigned via .Value
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}

Write-Host "Updated cryptos list data"